$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "solicitation" table (rows 34-41) ---

# Row 34: section header (merged E34:J34), styled like the other section
# headers (copy format from the "support" header row, E27:J27). Merge
# BEFORE pasting the format so the border treatment lands uniformly on
# every cell instead of only the merged range's outer edge.
$ws.Range("E34:J34").Merge()
$ws.Range("E27:J27").Copy()
$ws.Range("E34:J34").PasteSpecial(-4122)
$ws.Range("E34").Value = "solicitation"

# Row 35: column headers, styled like the other column-header rows
# (copy format from row 28).
$ws.Range("E28:J28").Copy()
$ws.Range("E35:J35").PasteSpecial(-4122)
$ws.Range("E35").Value = "ATRIBUTO"
$ws.Range("F35").Value = "Tipo de dados"
$ws.Range("G35").Value = "Comprimento"
$ws.Range("H35").Value = "Restrição"
$ws.Range("I35").Value = "PK"
$ws.Range("J35").Value = "FK"

# Rows 36-41: data rows, styled like the other data rows
# (copy format from row 7, which uses the plain centered/bordered style
# across every column).
$ws.Range("E7:J7").Copy()
$ws.Range("E36:J41").PasteSpecial(-4122)

$ws.Range("E36").Value = "ID"
$ws.Range("F36").Value = "String"
$ws.Range("G36").Value = "max"
$ws.Range("H36").Value = "not null"
$ws.Range("I36").Value = "sim"
$ws.Range("J36").Value = "—"

$ws.Range("E37").Value = "description"
$ws.Range("F37").Value = "String"
$ws.Range("G37").Value = "max"
$ws.Range("H37").Value = "not null"
$ws.Range("I37").Value = " —"
$ws.Range("J37").Value = "—"

$ws.Range("E38").Value = "approved"
$ws.Range("F38").Value = "Boolean"
$ws.Range("G38").Value = "max"
$ws.Range("H38").Value = "not null"
$ws.Range("I38").Value = "—"
$ws.Range("J38").Value = "—"

$ws.Range("E39").Value = "user_id"
$ws.Range("F39").Value = "String"
$ws.Range("G39").Value = "max"
$ws.Range("H39").Value = "not null"
$ws.Range("I39").Value = " —"
$ws.Range("J39").Value = "sim"

$ws.Range("E40").Value = "hospital_id"
$ws.Range("F40").Value = "String"
$ws.Range("G40").Value = "max"
$ws.Range("H40").Value = "not null"
$ws.Range("I40").Value = "—"
$ws.Range("J40").Value = "sim"

$ws.Range("E41").Value = "timestamps"
$ws.Range("F41").Value = "Date"
$ws.Range("G41").Value = "max"
$ws.Range("H41").Value = "not null"
$ws.Range("I41").Value = "—"
$ws.Range("J41").Value = "—"

# Match the selection left by the author after adding the table.
$ws.Range("E34:J41").Select()
